$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "670×8=" "217×7="
Replace-Text "206×5=" "682×8="
Replace-Text "427×7=" "239×8="
Replace-Text "847×2=" "654×5="
Replace-Text "953×3=" "460×7="
Replace-Text "519×5=" "279×4="
Replace-Text "265×9=" "210×6="
Replace-Text "843×9=" "384×3="
Replace-Text "155×3=" "367×6="
Replace-Text "365×3=" "366×4="
Replace-Text "657×6=" "157×7="
Replace-Text "626×4=" "674×6="
Replace-Text "931×6=" "386×9="
Replace-Text "970×6=" "351×5="
Replace-Text "585×3=" "761×3="
Replace-Text "113×9=" "740×6="
Replace-Text "443×6=" "159×6="
Replace-Text "321×6=" "842×5="
Replace-Text "187×3=" "818×3="
Replace-Text "799×7=" "968×5="
Replace-Text "198×8=" "879×2="
Replace-Text "314×8=" "466×2="
Replace-Text "971×5=" "214×6="
Replace-Text "588×5=" "180×3="
Replace-Text "866×6=" "205×7="
